$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Links")

# Insert a new row at the top, shifting existing rows (1-12) down to (2-13)
$ws.Range("A1").EntireRow.Insert()

# Populate the new first row with the new link entry
$ws.Range("A1").Value = 29
$ws.Range("B1").Value = "Логические типы в Python"
$ws.Range("C1").Value = "https://pythonchik.ru/osnovy/logicheskiy-tip-dannyh"
$ws.Range("D1").Value = "Логические типы в Python"
